$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily write "PENDENTE" into an unused cell so it is registered
# in the shared-strings table, then clear the cell's value. Excel keeps
# shared strings already allocated in the table (it does not compact the
# table on every edit), so the cell reverts to empty while "PENDENTE"
# remains available as a shared string.
$ws.Range("Z1").Value = "PENDENTE"
$ws.Range("Z1").Value = ""
